$d = $word.ActiveDocument

function Set-ParagraphInnerXml($doc, [string]$anchorText, [string]$innerRunXml) {
    $f = $doc.Content.Find
    $f.ClearFormatting()
    $f.Text = $anchorText
    $f.Forward = $true
    $f.Wrap = 1
    $found = $f.Execute()
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    $matchRange = $f.Parent
    $targetParagraph = $matchRange.Paragraphs(1)
    $paragraphRange = $targetParagraph.Range

    $payload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $innerRunXml + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    [void]$paragraphRange.InsertXML($payload)
}

# 1) "Programa" section: split the long run into multiple <w:t>/<w:br/> pieces.
$programaAnchor = "1.Solidificação / Fundição:Solidificação dos metais."
$programaRunXml = '<w:t>1.Solidificação / Fundição:</w:t><w:br/><w:t xml:space="preserve">Solidificação dos metais. Processos de montagem, tecnologia de fundição. Segregação e defeitos em peças fundidas. Qualidade de peças fundidas. Verificação experimental do efeito de algumas variáveis básicas de fundição tais como tipo de molde (cerâmicos, metálicos, ...), refrigeração ou pré-aquecimento do molde, temperatura de vazamento, adição de agentes inoculantes, agitação mecânica, dentre outras, nas micro e macroestruturas de metais e ligas metálicas. </w:t><w:br/><w:t>2.Experimentos de modificação superficial tais como: eletrodeposição, nitretação, carbonitretação, cementação, deposição por soldagem, aspersão Térmica, PVD, CVD, pack cementation, têmpera superficial, ensaios tribológicos (abrasão a dois ou mais corpos e roda de borracha), caracterizações microestruturais.</w:t><w:br/><w:br/><w:t>A turma será dividida em grupos, cada um recebendo orientação para realização de tarefas específicas. Dependendo do tipo de experimento, os resultados de todos os grupos serão compartilhados e cada grupo deverá desenvolver sua interpretação sobre o conjunto daqueles resultados.</w:t><w:br/><w:br/><w:t>Cada grupo deverá também desenvolver um trabalho (não experimental), em tema livre, procurando integrar temáticas de fundição/solidificação e modificação superficial.</w:t>'
Set-ParagraphInnerXml $d $programaAnchor $programaRunXml

# 2) "Bibliografia" section: split the long run into multiple <w:t>/<w:br/> pieces.
$biblioAnchor = "GARCIA, A. Solidificação: Fundamentos e Aplicações, Editora da Unicamp, 2001."
$biblioRunXml = '<w:t>GARCIA, A. Solidificação: Fundamentos e Aplicações, Editora da Unicamp, 2001.</w:t><w:br/><w:t>FLEMINGS, M.C. Solidification Processing, McGraw-Hill, 1974.</w:t><w:br/><w:t>CAMPOS FILHO, M.P. de, DAVIES, G.J. - Solidificação e fundição de metais e suas ligas. S.Paulo, EDUSP, 1978.</w:t><w:br/><w:t>SIEGEL, M. - Fundição. S.Paulo, ABM, 1979.</w:t><w:br/><w:t>Krauss, G. "Steels: Heat Treating and Processing Principles" - ASM International, 1990.</w:t><w:br/><w:t>ASM Handbook, vol. 5 – Surface Engineering. ISBN: 978-0-87170-384-2</w:t><w:br/><w:t>PLETCHER, D.; WALSH, F.C. Industrial Electrochemistry, 2ª Ed. Springer, 1990.</w:t>'
Set-ParagraphInnerXml $d $biblioAnchor $biblioRunXml

Write-Host "Done."
